# Update the "想去人数" (F column) counts on several rows across sheets
# to reflect a refreshed data pull (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1356
$ws1.Range("F6").Value = 40562
$ws1.Range("F11").Value = 5563
$ws1.Range("F12").Value = 385
$ws1.Range("F14").Value = 2663
$ws1.Range("F15").Value = 6252
$ws1.Range("F16").Value = 161
$ws1.Range("F17").Value = 1162
$ws1.Range("F20").Value = 5
$ws1.Range("F34").Value = 173
$ws1.Range("F36").Value = 183
$ws1.Range("F38").Value = 34

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F28").Value = 385
$ws2.Range("F38").Value = 41

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1356
$ws4.Range("F17").Value = 5563
$ws4.Range("F18").Value = 385
$ws4.Range("F20").Value = 2663
$ws4.Range("F22").Value = 6252
$ws4.Range("F24").Value = 161
$ws4.Range("F25").Value = 1162
$ws4.Range("F42").Value = 173
$ws4.Range("F44").Value = 183
$ws4.Range("F50").Value = 41
